$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update bus voltage magnitude results (res_bus/vm_pu) for the "380 kV" case.
# Slack/reference bus setpoint moves from 1.05 to 1.02 p.u. (column B), which
# changes the power-flow solution for every other bus (columns C:F, I:N).
# Columns G (slack = 1) and H (unused) are unchanged.

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.042787090767824
$ws.Cells.Item(2, 4).Value = 1.046997488767504
$ws.Cells.Item(2, 5).Value = 1.040796569429764
$ws.Cells.Item(2, 6).Value = 1.055794405078732
$ws.Cells.Item(2, 9).Value = 1.038020593043764
$ws.Cells.Item(2, 10).Value = 1.047861036300901
$ws.Cells.Item(2, 11).Value = 1.049761358714254
$ws.Cells.Item(2, 12).Value = 1.043577901587617
$ws.Cells.Item(2, 13).Value = 1.058533883661371
$ws.Cells.Item(2, 14).Value = 1.049349119163789

$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.044370604914972
$ws.Cells.Item(3, 4).Value = 1.048217761824661
$ws.Cells.Item(3, 5).Value = 1.042165840438479
$ws.Cells.Item(3, 6).Value = 1.05725359839982
$ws.Cells.Item(3, 9).Value = 1.038419654363171
$ws.Cells.Item(3, 10).Value = 1.049087984523825
$ws.Cells.Item(3, 11).Value = 1.050792327144732
$ws.Cells.Item(3, 12).Value = 1.044756200884216
$ws.Cells.Item(3, 13).Value = 1.059804943729312
$ws.Cells.Item(3, 14).Value = 1.050577809793923

$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.045393329847113
$ws.Cells.Item(4, 4).Value = 1.049005464576918
$ws.Cells.Item(4, 5).Value = 1.043050311366623
$ws.Cells.Item(4, 6).Value = 1.058196268679318
$ws.Cells.Item(4, 9).Value = 1.038675331154419
$ws.Cells.Item(4, 10).Value = 1.049879646213311
$ws.Cells.Item(4, 11).Value = 1.051456978670331
$ws.Cells.Item(4, 12).Value = 1.045516576333017
$ws.Cells.Item(4, 13).Value = 1.060625367522939
$ws.Cells.Item(4, 14).Value = 1.051370595733813

$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.045822834730403
$ws.Cells.Item(5, 4).Value = 1.049336167425704
$ws.Cells.Item(5, 5).Value = 1.043421781993641
$ws.Cells.Item(5, 6).Value = 1.058592209718227
$ws.Cells.Item(5, 9).Value = 1.038782211709885
$ws.Cells.Item(5, 10).Value = 1.050211927988435
$ws.Cells.Item(5, 11).Value = 1.051735816658601
$ws.Cells.Item(5, 12).Value = 1.045835751171845
$ws.Cells.Item(5, 13).Value = 1.060969792978182
$ws.Cells.Item(5, 14).Value = 1.05170334938717

$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.045894924459653
$ws.Cells.Item(6, 4).Value = 1.049391667839162
$ws.Cells.Item(6, 5).Value = 1.043484132596577
$ws.Cells.Item(6, 6).Value = 1.05865866921486
$ws.Cells.Item(6, 9).Value = 1.038800121969972
$ws.Cells.Item(6, 10).Value = 1.050267688526854
$ws.Cells.Item(6, 11).Value = 1.051782600856321
$ws.Cells.Item(6, 12).Value = 1.045889313672917
$ws.Cells.Item(6, 13).Value = 1.061027595540253
$ws.Cells.Item(6, 14).Value = 1.05175918911195

$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.045399070663446
$ws.Cells.Item(7, 4).Value = 1.049009885194629
$ws.Cells.Item(7, 5).Value = 1.043055276384388
$ws.Cells.Item(7, 6).Value = 1.058201560658371
$ws.Cells.Item(7, 9).Value = 1.038676761675338
$ws.Cells.Item(7, 10).Value = 1.049884088263262
$ws.Cells.Item(7, 11).Value = 1.051460706794789
$ws.Cells.Item(7, 12).Value = 1.045520843067937
$ws.Cells.Item(7, 13).Value = 1.060629971633009
$ws.Cells.Item(7, 14).Value = 1.051375044091984

$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.043322649066866
$ws.Cells.Item(8, 4).Value = 1.047410282106977
$ws.Cells.Item(8, 5).Value = 1.041259644151801
$ws.Cells.Item(8, 6).Value = 1.056287866741372
$ws.Cells.Item(8, 9).Value = 1.038155985790313
$ws.Cells.Item(8, 10).Value = 1.048276160670448
$ws.Cells.Item(8, 11).Value = 1.050110291048095
$ws.Cells.Item(8, 12).Value = 1.043976544663898
$ws.Cells.Item(8, 13).Value = 1.058963870360461
$ws.Cells.Item(8, 14).Value = 1.049764833057556

$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.039648629710175
$ws.Cells.Item(9, 4).Value = 1.044576771155948
$ws.Cells.Item(9, 5).Value = 1.038083387239365
$ws.Cells.Item(9, 6).Value = 1.052903667343646
$ws.Cells.Item(9, 9).Value = 1.037218715799758
$ws.Cells.Item(9, 10).Value = 1.045425182862012
$ws.Cells.Item(9, 11).Value = 1.047711621922354
$ws.Cells.Item(9, 12).Value = 1.041239183684626
$ws.Cells.Item(9, 13).Value = 1.056012071558581
$ws.Cells.Item(9, 14).Value = 1.046909806533621

$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.037188519987097
$ws.Cells.Item(10, 4).Value = 1.042677402555624
$ws.Cells.Item(10, 5).Value = 1.035957257988397
$ws.Cells.Item(10, 6).Value = 1.050638965035731
$ws.Cells.Item(10, 9).Value = 1.036580519193182
$ws.Cells.Item(10, 10).Value = 1.043512235418447
$ws.Cells.Item(10, 11).Value = 1.046099311983451
$ws.Cells.Item(10, 12).Value = 1.039403009281908
$ws.Cells.Item(10, 13).Value = 1.054033062926468
$ws.Cells.Item(10, 14).Value = 1.044994142485268

$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.036120575533951
$ws.Cells.Item(11, 4).Value = 1.041852406017325
$ws.Cells.Item(11, 5).Value = 1.035034468790766
$ws.Cells.Item(11, 6).Value = 1.04965618287657
$ws.Cells.Item(11, 9).Value = 1.036300967144768
$ws.Cells.Item(11, 10).Value = 1.042680888630579
$ws.Cells.Item(11, 11).Value = 1.045397948463861
$ws.Cells.Item(11, 12).Value = 1.038605155425506
$ws.Cells.Item(11, 13).Value = 1.053173388166786
$ws.Cells.Item(11, 14).Value = 1.044161615089604

$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.03572347579017
$ws.Cells.Item(12, 4).Value = 1.04154557337099
$ws.Cells.Item(12, 5).Value = 1.034691369872101
$ws.Cells.Item(12, 6).Value = 1.049290801196219
$ws.Cells.Item(12, 9).Value = 1.036196643513766
$ws.Cells.Item(12, 10).Value = 1.042371625530534
$ws.Cells.Item(12, 11).Value = 1.045136939200485
$ws.Cells.Item(12, 12).Value = 1.0383083710803
$ws.Cells.Item(12, 13).Value = 1.052853644729448
$ws.Cells.Item(12, 14).Value = 1.043851912800483

$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.035808674123303
$ws.Cells.Item(13, 4).Value = 1.04161140789197
$ws.Cells.Item(13, 5).Value = 1.034764981014617
$ws.Cells.Item(13, 6).Value = 1.049369192011728
$ws.Cells.Item(13, 9).Value = 1.036219043315352
$ws.Cells.Item(13, 10).Value = 1.042437984652203
$ws.Cells.Item(13, 11).Value = 1.045192948940963
$ws.Cells.Item(13, 12).Value = 1.038372051739936
$ws.Cells.Item(13, 13).Value = 1.052922250016307
$ws.Cells.Item(13, 14).Value = 1.043918366159717

$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.036087759733711
$ws.Cells.Item(14, 4).Value = 1.04182705118224
$ws.Cells.Item(14, 5).Value = 1.035006114992735
$ws.Cells.Item(14, 6).Value = 1.049625987153372
$ws.Cells.Item(14, 9).Value = 1.036292353646298
$ws.Cells.Item(14, 10).Value = 1.042655334361642
$ws.Cells.Item(14, 11).Value = 1.045376383429087
$ws.Cells.Item(14, 12).Value = 1.038580631876608
$ws.Cells.Item(14, 13).Value = 1.053146966732612
$ws.Cells.Item(14, 14).Value = 1.044136024530673

$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.036259657893075
$ws.Cells.Item(15, 4).Value = 1.041959863920282
$ws.Cells.Item(15, 5).Value = 1.035154641069495
$ws.Cells.Item(15, 6).Value = 1.049784162754843
$ws.Cells.Item(15, 9).Value = 1.03633745811543
$ws.Cells.Item(15, 10).Value = 1.042789188951522
$ws.Cells.Item(15, 11).Value = 1.045489338102138
$ws.Cells.Item(15, 12).Value = 1.03870908829147
$ws.Cells.Item(15, 13).Value = 1.053285365944734
$ws.Cells.Item(15, 14).Value = 1.044270069209425

$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.037259338031462
$ws.Cells.Item(16, 4).Value = 1.042732100298116
$ws.Cells.Item(16, 5).Value = 1.036018454148585
$ws.Cells.Item(16, 6).Value = 1.050704142925058
$ws.Cells.Item(16, 9).Value = 1.036599004241923
$ws.Cells.Item(16, 10).Value = 1.043567344607699
$ws.Cells.Item(16, 11).Value = 1.046145790637772
$ws.Cells.Item(16, 12).Value = 1.039455900956809
$ws.Cells.Item(16, 13).Value = 1.05409005801265
$ws.Cells.Item(16, 14).Value = 1.045049329935888

$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.037885679565925
$ws.Cells.Item(17, 4).Value = 1.043215813478868
$ws.Cells.Item(17, 5).Value = 1.036559716168132
$ws.Cells.Item(17, 6).Value = 1.05128063962384
$ws.Cells.Item(17, 9).Value = 1.036762203709467
$ws.Cells.Item(17, 10).Value = 1.044054644120721
$ws.Cells.Item(17, 11).Value = 1.046556697929132
$ws.Cells.Item(17, 12).Value = 1.039923607023969
$ws.Cells.Item(17, 13).Value = 1.054594077492998
$ws.Cells.Item(17, 14).Value = 1.045537321470114

$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.038250754214631
$ws.Cells.Item(18, 4).Value = 1.043497708984961
$ws.Cells.Item(18, 5).Value = 1.036875217250142
$ws.Cells.Item(18, 6).Value = 1.051616693486598
$ws.Cells.Item(18, 9).Value = 1.036857085841639
$ws.Cells.Item(18, 10).Value = 1.044338586142813
$ws.Cells.Item(18, 11).Value = 1.046796062766384
$ws.Cells.Item(18, 12).Value = 1.040196144580417
$ws.Cells.Item(18, 13).Value = 1.054887798689914
$ws.Cells.Item(18, 14).Value = 1.045821666722446

$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.038375191531118
$ws.Cells.Item(19, 4).Value = 1.043593786546406
$ws.Cells.Item(19, 5).Value = 1.036982759980638
$ws.Cells.Item(19, 6).Value = 1.05173124435839
$ws.Cells.Item(19, 9).Value = 1.036889385829161
$ws.Cells.Item(19, 10).Value = 1.044435353854573
$ws.Cells.Item(19, 11).Value = 1.04687762755934
$ws.Cells.Item(19, 12).Value = 1.04028902773083
$ws.Cells.Item(19, 13).Value = 1.054987905366697
$ws.Cells.Item(19, 14).Value = 1.045918571855456

$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.037818506002195
$ws.Cells.Item(20, 4).Value = 1.043163941156641
$ws.Cells.Item(20, 5).Value = 1.036501665446397
$ws.Cells.Item(20, 6).Value = 1.051218808453832
$ws.Cells.Item(20, 9).Value = 1.036744725974457
$ws.Cells.Item(20, 10).Value = 1.044002391714979
$ws.Cells.Item(20, 11).Value = 1.046512643629125
$ws.Cells.Item(20, 12).Value = 1.039873454317491
$ws.Cells.Item(20, 13).Value = 1.054540028428889
$ws.Cells.Item(20, 14).Value = 1.045484994859962

$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.03600558760181
$ws.Cells.Item(21, 4).Value = 1.041763560499107
$ws.Cells.Item(21, 5).Value = 1.034935116338732
$ws.Cells.Item(21, 6).Value = 1.04955037664309
$ws.Cells.Item(21, 9).Value = 1.036270779012839
$ws.Cells.Item(21, 10).Value = 1.04259134316758
$ws.Cells.Item(21, 11).Value = 1.045322380175772
$ws.Cells.Item(21, 12).Value = 1.03851922205063
$ws.Cells.Item(21, 13).Value = 1.053080804978501
$ws.Cells.Item(21, 14).Value = 1.044071942461776

$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.034863311649038
$ws.Cells.Item(22, 4).Value = 1.040880812446558
$ws.Cells.Item(22, 5).Value = 1.033948228128765
$ws.Cells.Item(22, 6).Value = 1.048499438115818
$ws.Cells.Item(22, 9).Value = 1.035969978734079
$ws.Cells.Item(22, 10).Value = 1.041701472525164
$ws.Cells.Item(22, 11).Value = 1.044571166171135
$ws.Cells.Item(22, 12).Value = 1.037665294090193
$ws.Cells.Item(22, 13).Value = 1.05216088763296
$ws.Cells.Item(22, 14).Value = 1.043180808100979

$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.03546908733163
$ws.Cells.Item(23, 4).Value = 1.041348992043643
$ws.Cells.Item(23, 5).Value = 1.034471582911527
$ws.Cells.Item(23, 6).Value = 1.049056746558793
$ws.Cells.Item(23, 9).Value = 1.036129706276681
$ws.Cells.Item(23, 10).Value = 1.042173467638587
$ws.Cells.Item(23, 11).Value = 1.044969671405235
$ws.Cells.Item(23, 12).Value = 1.038118214291315
$ws.Cells.Item(23, 13).Value = 1.052648787952345
$ws.Cells.Item(23, 14).Value = 1.043653473501601

$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.037848859646957
$ws.Cells.Item(24, 4).Value = 1.043187380787181
$ws.Cells.Item(24, 5).Value = 1.036527896711321
$ws.Cells.Item(24, 6).Value = 1.051246747937584
$ws.Cells.Item(24, 9).Value = 1.036752624366554
$ws.Cells.Item(24, 10).Value = 1.044026003228911
$ws.Cells.Item(24, 11).Value = 1.046532550830354
$ws.Cells.Item(24, 12).Value = 1.039896116993363
$ws.Cells.Item(24, 13).Value = 1.054564451692938
$ws.Cells.Item(24, 14).Value = 1.045508639904953

$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.040600303766098
$ws.Cells.Item(25, 4).Value = 1.045311097855878
$ws.Cells.Item(25, 5).Value = 1.03890601097843
$ws.Cells.Item(25, 6).Value = 1.053780038695306
$ws.Cells.Item(25, 9).Value = 1.037463362666434
$ws.Cells.Item(25, 10).Value = 1.046164362191418
$ws.Cells.Item(25, 11).Value = 1.048334033428975
$ws.Cells.Item(25, 12).Value = 1.041948810839922
$ws.Cells.Item(25, 13).Value = 1.056777113631101
$ws.Cells.Item(25, 14).Value = 1.047650035582459
